$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.023.26'
$ws.Range('E2').Value = '  -1.20%  '

$ws.Range('D3').Value = '2.676.92'
$ws.Range('E3').Value = '  +0.80%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '''599.60'
$ws.Range('E5').Value = '  -1.25%  '

$ws.Range('D6').Value = '''175.46'
$ws.Range('E6').Value = '  -3.05%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('E8').Value = '  -0.94%  '

$ws.Range('D9').Value = '2.675.09'
$ws.Range('E9').Value = '  +0.77%  '

$ws.Range('E10').Value = '  -2.98%  '

$ws.Range('E11').Value = '  +2.17%  '

$ws.Range('D12').Value = '''0.357'
$ws.Range('E12').Value = '  +1.24%  '

$ws.Range('D13').Value = '''4.99'

$ws.Range('D14').Value = '3.171.93'
$ws.Range('E14').Value = '  +1.18%  '

$ws.Range('E15').Value = '  -3.37%  '

$ws.Range('D16').Value = '71.836.71'
$ws.Range('E16').Value = '  -1.24%  '

$ws.Range('D17').Value = '''26.27'
$ws.Range('E17').Value = '  -2.71%  '

$ws.Range('D18').Value = '2.679.15'
$ws.Range('E18').Value = '  +0.64%  '

$ws.Range('D19').Value = '''12.26'
$ws.Range('E19').Value = '  +5.71%  '

$ws.Range('D20').Value = '''8.22'
$ws.Range('E20').Value = '  +3.66%  '

$ws.Range('D21').Value = '''372.56'
$ws.Range('E21').Value = '  -2.77%  '

$ws.Range('E22').Value = '  -1.52%  '

$ws.Range('D23').Value = '''2.03'
$ws.Range('E23').Value = '  +0.01%  '

$ws.Range('D24').Value = '''72.04'
$ws.Range('E24').Value = '  -2.02%  '

$ws.Range('D25').Value = '''1.00'
$ws.Range('E25').Value = '  -0.06%  '

$ws.Range('D26').Value = '''4.35'
$ws.Range('E26').Value = '  -2.19%  '

$ws.Range('D27').Value = '''9.80'
$ws.Range('E27').Value = '  -2.00%  '

$ws.Range('D28').Value = '2.817.36'
$ws.Range('E28').Value = '  +0.90%  '

$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.31%  '

$ws.Range('D30').Value = '0.0₃0974'
$ws.Range('E30').Value = '  -0.09%  '

$ws.Range('D31').Value = '''8.08'
$ws.Range('E31').Value = '  -0.37%  '

$ws.Range('D32').Value = '''502.40'
$ws.Range('E32').Value = '  -7.26%  '

$ws.Range('D33').Value = '''1.31'
$ws.Range('E33').Value = '  -3.01%  '

$ws.Range('D34').Value = '''1.83'
$ws.Range('E34').Value = '  -1.18%  '

$ws.Range('E35').Value = '  +0.01%  '

$ws.Range('D36').Value = '''162.32'
$ws.Range('E36').Value = '  -0.09%  '

$ws.Range('D37').Value = '''19.59'
$ws.Range('E37').Value = '  +1.01%  '

$ws.Range('D38').Value = '''19.07'
$ws.Range('E38').Value = '  -0.28%  '

$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '''0.110'
$ws.Range('E39').Value = '  -2.38%  '

$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').Value = '''1.38'
$ws.Range('E40').Value = '  -2.84%  '

$ws.Range('D41').Value = '''1.78'
$ws.Range('E41').Value = '  -4.43%  '

$ws.Range('E42').Value = '  -0.10%  '

$ws.Range('D43').Value = '''5.01'

$ws.Range('E44').Value = '  -3.50%  '

$ws.Range('E45').Value = '  -0.78%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''156.43'
$ws.Range('E46').Value = '  +2.95%  '

$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '''39.49'
$ws.Range('E47').Value = '  -0.75%  '

$ws.Range('E48').Value = '  +2.80%  '

$ws.Range('D49').Value = '''3.73'
$ws.Range('E49').Value = '  +0.63%  '

$ws.Range('E50').Value = '  +1.50%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.0755'
$ws.Range('E51').Value = '  -1.75%  '

